$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..39 already contain Data / Entrada1 / Saida1 in columns A-C,
# with column D ("Entrada2") always holding the placeholder "S".
# Add more time-card detail: give Entrada2 a real value, and add
# Saida2 / Entrada3 columns (E, F) for each record row.
$lastRow = 39

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "13:45"
    $ws.Cells.Item($r, 5).Value = "14:00"
    $ws.Cells.Item($r, 6).Value = "S"
}
